$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# --- Update the Right/Wrong/Not-Attempt summary (row 10) and Totals (row 12) ---
# The student did not attempt the quiz: 0 right, 0 wrong, all 28 not attempted -> Absent.
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 28

$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = "Absent"

# --- Clear the student's submitted answers (column A rows 16-40, column D rows 16-18) ---
# The "Correct Ans" columns (B and E) are left untouched.
# A24 already carries the plain "normalStyle" (no highlight) format used for a
# blank "Student Ans" cell, so copy that format onto the cells being cleared.
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A16:A40").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D16:D18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A16:A40").ClearContents() | Out-Null
$ws.Range("D16:D18").ClearContents() | Out-Null
